# Generate Report for Handoff
#
# The a9abe794-... file has moved from "Handed back: in sync with en-US"
# to "Ready for handoff" on every sheet (Overview + per-locale sheets),
# and the "Latest Handoff Datetime" for that run was refreshed for both
# locales (zh-cn: 17:20:27 -> 17:22:07, de-de: 17:20:43 -> 17:22:15).

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"

# --- Overview sheet: a9abe794 row (row 3) status for both locale columns ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $status
$overview.Range("C3").Value = $status
# Overview "Latest Handoff Date" column reflects the newest handoff timestamp
$overview.Range("D2").Value = "2016-03-19 17:22:15"
$overview.Range("D3").Value = "2016-03-19 17:22:15"

# --- zh-cn sheet: a9abe794 row (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $status
$zhcn.Range("E2").Value = "2016-03-19 17:22:07"
$zhcn.Range("E3").Value = "2016-03-19 17:22:07"

# --- de-de sheet: a9abe794 row (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $status
$dede.Range("E2").Value = "2016-03-19 17:22:15"
$dede.Range("E3").Value = "2016-03-19 17:22:15"
